$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 74.57895000000001
$ws.Range("I15").Value = 74.57895000000001
$ws.Range("K15").Value = 223.73685
$ws.Range("M15").Value = -54.73685

$ws.Range("H51").Value = 11388.444
$ws.Range("I51").Value = 9582.833000000001
$ws.Range("K51").Value = 9582.833000000001
$ws.Range("M51").Value = -9098.833000000001

$ws.Range("H113").Value = 23284.428
$ws.Range("I113").Value = 25498.5
$ws.Range("J113").Value = 10000
$ws.Range("K113").Value = 25498.5
$ws.Range("L113").Value = 10000
$ws.Range("M113").Value = -22244.5
$ws.Range("N113").Value = -16508

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 89
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 89
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 89
$ws.Range("N4").Value = -321
$ws.Range("M4").ClearContents()

$ws.Range("H32").Value = 3911.8
$ws.Range("I32").Value = 3542.818
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 3542.818
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -3255.818
$ws.Range("N32").Value = -10574

$ws.Range("H45").Value = 4617.4287
$ws.Range("I45").Value = 3915
$ws.Range("K45").Value = 3915
$ws.Range("M45").Value = -3538

$ws.Range("H61").Value = 2249.5
$ws.Range("I61").Value = 1499
$ws.Range("K61").Value = 1499
$ws.Range("M61").Value = -1287

$ws.Range("H62").Value = 50000
$ws.Range("J62").Value = 50000
$ws.Range("L62").Value = 50000
$ws.Range("N62").Value = -51248

$ws.Range("H65").Value = 50000
$ws.Range("J65").Value = 50000
$ws.Range("L65").Value = 150000
$ws.Range("N65").Value = -156240

$ws.Range("H97").Value = 1512.75
$ws.Range("I97").Value = 1410.5
$ws.Range("K97").Value = 1410.5
$ws.Range("M97").Value = -914.5

$ws.Range("H110").Value = 10000
$ws.Range("I110").Value = 10000
$ws.Range("K110").Value = 10000
$ws.Range("M110").Value = -7955

$ws.Range("H136").Value = 2249.5
$ws.Range("I136").Value = 1499
$ws.Range("K136").Value = 4497
$ws.Range("M136").Value = -1947

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 18750
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 18750
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 18750
$ws.Range("N3").Value = -18976
$ws.Range("M3").ClearContents()

$ws.Range("H4").Value = 19987.5
$ws.Range("I4").Value = 4150
$ws.Range("J4").Value = 27906.25
$ws.Range("K4").Value = 4150
$ws.Range("L4").Value = 27906.25
$ws.Range("N4").Value = -28130.25
$ws.Range("M4").Value = -4038

$ws.Range("H11").Value = 9550
$ws.Range("I11").Value = 750
$ws.Range("J11").Value = 13950
$ws.Range("K11").Value = 750
$ws.Range("L11").Value = 13950
$ws.Range("M11").Value = -610
$ws.Range("N11").Value = -14230

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H112").Value = 29000
$ws.Range("I112").Value = 28000
$ws.Range("K112").Value = 84000
$ws.Range("M112").Value = -82892

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5000000
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").ClearContents()

$ws.Range("H21").Value = 39002.332
$ws.Range("J21").Value = 39002.332
$ws.Range("L21").Value = 39002.332
$ws.Range("N21").Value = -39348.332

$ws.Range("H30").Value = 39002.332
$ws.Range("J30").Value = 39002.332
$ws.Range("L30").Value = 39002.332
$ws.Range("N30").Value = -39212.332

$ws.Range("H92").Value = 13562.75
$ws.Range("J92").Value = 13562.75
$ws.Range("L92").Value = 13562.75
$ws.Range("N92").Value = -17306.75

$ws.Range("H122").Value = 3705.7273
$ws.Range("I122").Value = 3668.111
$ws.Range("K122").Value = 11004.333
$ws.Range("M122").Value = -8554.332999999999

$ws.Range("H126").Value = 7509
$ws.Range("I126").Value = 6610.8
$ws.Range("K126").Value = 19832.4
$ws.Range("M126").Value = -17362.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3258.4
$ws.Range("I22").Value = 3258.4
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 3258.4
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -2963.4
$ws.Range("N22").ClearContents()

$ws.Range("H23").Value = 18000
$ws.Range("I23").Value = 18000
$ws.Range("K23").Value = 18000
$ws.Range("M23").Value = -17770

$ws.Range("H27").Value = 3258.4
$ws.Range("I27").Value = 3258.4
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 3258.4
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -3151.4
$ws.Range("N27").ClearContents()

$ws.Range("H31").Value = 3000
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3000
$ws.Range("N31").Value = -3496
$ws.Range("M31").ClearContents()

$ws.Range("I82").Value = 3001
$ws.Range("J82").Value = 3700
$ws.Range("K82").Value = 3001
$ws.Range("L82").Value = 3700
$ws.Range("M82").Value = -2640
$ws.Range("N82").Value = -4422

$ws.Range("I85").Value = 3001
$ws.Range("J85").Value = 3700
$ws.Range("K85").Value = 3001
$ws.Range("L85").Value = 3700
$ws.Range("M85").Value = -1753
$ws.Range("N85").Value = -6196

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 2333.3333
$ws.Range("J5").Value = 2333.3333
$ws.Range("L5").Value = 2333.3333
$ws.Range("N5").Value = -2557.3333

$ws.Range("H94").Value = 30329.5
$ws.Range("J94").Value = 30329.5
$ws.Range("L94").Value = 30329.5
$ws.Range("N94").Value = -32131.5

Write-Output "done"
